$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 33452.02960683603
$ws.Range("B3").Value = 33416.04744270352
$ws.Range("B4").Value = 33428.64858119172
$ws.Range("B5").Value = 33334.21667591591
$ws.Range("B6").Value = 33184.44809034982
$ws.Range("B7").Value = 32975.84078225301
$ws.Range("B8").Value = 33011.55169400804
$ws.Range("B9").Value = 32975.13320882477
$ws.Range("B10").Value = 32805.19348058276
$ws.Range("B11").Value = 32634.48512479875
$ws.Range("B12").Value = 32508.79416243713
$ws.Range("B13").Value = 32482.02259935138
$ws.Range("B14").Value = 32510.15110191844
$ws.Range("B15").Value = 32309.01164024487
$ws.Range("B16").Value = 32303.57034460452
$ws.Range("B17").Value = 32368.46900602943
$ws.Range("B18").Value = 32424.29800875598
$ws.Range("B19").Value = 32289.62553800494
$ws.Range("B20").Value = 32200.02935571772
$ws.Range("B21").Value = 32242.01235459997
$ws.Range("B22").Value = 32351.16323853311
$ws.Range("B23").Value = 32105.96981491959
$ws.Range("B24").Value = 31943.94469998771
$ws.Range("B25").Value = 32073.11814260486
$ws.Range("B26").Value = 32154.88694852356
$ws.Range("B27").Value = 32208.2683725684
$ws.Range("B28").Value = 31908.54353754226
$ws.Range("B29").Value = 30642.85604881407
$ws.Range("B30").Value = 31068.17665796504
$ws.Range("B31").Value = 31008.83988945493
$ws.Range("B32").Value = 30734.55695843923
$ws.Range("B33").Value = 30611.35030012325
$ws.Range("B34").Value = 30799.50791153299
$ws.Range("B35").Value = 30967.59615007453
$ws.Range("B36").Value = 30873.29989747982
$ws.Range("B37").Value = 30796.53690408655
$ws.Range("B38").Value = 30551.13515675055
$ws.Range("B39").Value = 30575.81542639952
$ws.Range("B40").Value = 30572.6086179704
$ws.Range("B41").Value = 30547.83795904304
$ws.Range("B42").Value = 30471.44441512788
$ws.Range("B43").Value = 30470.6776355125
$ws.Range("B44").Value = 30356.83820936358
$ws.Range("B45").Value = 30213.1546442889
$ws.Range("B46").Value = 30212.38023607982
$ws.Range("B47").Value = 30297.98927386082
$ws.Range("B48").Value = 30385.26081484321
$ws.Range("B49").Value = 30467.13777328087
$ws.Range("B50").Value = 30319.03494357815
$ws.Range("B51").Value = 29686.55449616881
$ws.Range("B52").Value = 29756.91414088169
$ws.Range("B53").Value = 29839.04640216875
$ws.Range("B54").Value = 29820.55082688458
$ws.Range("B55").Value = 29649.71324350654
$ws.Range("B56").Value = 29743.23005834155
$ws.Range("B57").Value = 29841.84302143726
$ws.Range("B58").Value = 29596.19076996803
$ws.Range("B59").Value = 29450.81339309886
$ws.Range("B60").Value = 29543.19129274989
$ws.Range("B61").Value = 29634.8434108577
$ws.Range("B62").Value = 29219.29997703242
$ws.Range("B63").Value = 28966.57561139015
$ws.Range("B64").Value = 29057.33278001472
$ws.Range("B65").Value = 29145.38188447343
$ws.Range("B66").Value = 29246.95598604753
$ws.Range("B67").Value = 28993.89488688272
$ws.Range("B68").Value = 28681.97790506395
$ws.Range("B69").Value = 28771.36731105049
$ws.Range("B70").Value = 28864.64162605524
$ws.Range("B71").Value = 28458.24801734745
$ws.Range("B72").Value = 28106.88352790985
$ws.Range("B73").Value = 28162.76178356502
$ws.Range("B74").Value = 28264.33588513912
$ws.Range("B75").Value = 27622.83676622821
$ws.Range("B76").Value = 27205.37591407311
$ws.Range("B77").Value = 27264.59249876102
$ws.Range("B78").Value = 27344.47754123756
$ws.Range("B79").Value = 27444.4684550722
$ws.Range("B80").Value = 27035.47954344309
$ws.Range("B81").Value = 26610.90267570942
$ws.Range("B82").Value = 26622.94960653624
$ws.Range("B83").Value = 26711.35061974639
$ws.Range("B84").Value = 26400.83758162509
$ws.Range("B85").Value = 25972.40364827805
$ws.Range("B86").Value = 26022.68150651307
$ws.Range("B87").Value = 26098.74984110688
$ws.Range("B88").Value = 25856.61943233096
$ws.Range("B89").Value = 25662.97805429606
$ws.Range("B90").Value = 25454.77825319822
$ws.Range("B91").Value = 25549.58160348714
$ws.Range("B92").Value = 25646.26474413522
$ws.Range("B93").Value = 25737.63981494185
$ws.Range("B94").Value = 25802.63411897263
$ws.Range("B95").Value = 25841.04890696665
$ws.Range("B96").Value = 25942.62300854076
$ws.Range("B97").Value = 25934.18438418673
$ws.Range("B98").Value = 25852.89319805034
$ws.Range("B99").Value = 25948.64422270136
$ws.Range("B100").Value = 26045.92514478829
$ws.Range("B101").Value = 26147.49924636239
$ws.Range("B102").Value = 26027.0704228528
$ws.Range("B103").Value = 26064.52119328914
$ws.Range("B104").Value = 26166.09529486324
$ws.Range("B105").Value = 26259.81010116188
$ws.Range("B106").Value = 26302.85011935838
$ws.Range("B107").Value = 26005.76017535929
$ws.Range("B108").Value = 26107.33427693339
$ws.Range("B109").Value = 26208.90837850749
$ws.Range("B110").Value = 26179.04562046386
$ws.Range("B111").Value = 25989.09310025768
$ws.Range("B112").Value = 26071.64056190167
$ws.Range("B113").Value = 26165.46444988354
$ws.Range("B114").Value = 26127.18647393826
$ws.Range("B115").Value = 26121.62757055038
$ws.Range("B116").Value = 26050.22183603783
$ws.Range("B117").Value = 26148.45198856734
$ws.Range("B118").Value = 26250.02609014144
$ws.Range("B119").Value = 26204.22692825878
$ws.Range("B120").Value = 26102.39756841367
$ws.Range("B121").Value = 26202.47713026175
$ws.Range("B122").Value = 26296.07940089858
$ws.Range("B123").Value = 26008.99881118023
$ws.Range("B124").Value = 25890.17297842748
$ws.Range("B125").Value = 25978.15484699187
$ws.Range("B126").Value = 26078.76704565335
$ws.Range("B127").Value = 26158.52771133001
$ws.Range("B128").Value = 25988.51922817296
$ws.Range("B129").Value = 24597.16328837206
$ws.Range("B130").Value = 24692.4545682901
$ws.Range("B131").Value = 24788.20116986421
$ws.Range("B132").Value = 24505.99523804776
$ws.Range("B133").Value = 23489.87557854199
$ws.Range("B134").Value = 23587.36316171896
$ws.Range("B135").Value = 23667.28162226742
$ws.Range("B136").Value = 23411.50229448806
$ws.Range("B137").Value = 23185.90689035285
$ws.Range("B138").Value = 23242.00039397286
$ws.Range("B139").Value = 23280.9014151942
$ws.Range("B140").Value = 23350.36680615055
$ws.Range("B141").Value = 23348.83670219993
$ws.Range("B142").Value = 23350.91560418215
$ws.Range("B143").Value = 23442.88318346325
$ws.Range("B144").Value = 23544.45728503736
$ws.Range("B145").Value = 23635.01414302171
$ws.Range("B146").Value = 23727.30991697343
$ws.Range("B147").Value = 23828.88401854753
$ws.Range("B148").Value = 23930.45812012163
$ws.Range("B149").Value = 23756.50599891226
$ws.Range("B150").Value = 23538.90477617816
$ws.Range("B151").Value = 23640.47887775226
